$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.210.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.825.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "453.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.739"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000317"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.87%  "
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.430.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.48%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.137"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.800.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("E19").Value = "  +8.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.288.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.67%  "
$ws.Range("E22").Value = "  +3.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +16.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "37.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "736.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.135"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.26%  "
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.10%  "
$ws.Range("E35").Value = "  +4.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.24%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0474"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.353"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.26%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0688"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.28%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +4.73%  "
$ws.Range("E46").Value = "  +4.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.69%  "
$ws.Range("E48").Value = "  +5.44%  "
$ws.Range("E49").Value = "  +4.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +2.96%  "
